$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 2128.3125
$ws.Range("J39").Value = 1368.1
$ws.Range("L39").Value = 4104.299999999999
$ws.Range("N39").Value = -4696.299999999999
$ws.Range("H54").Value = 7499.5
$ws.Range("I54").Value = 9999
$ws.Range("J54").Value = 5000
$ws.Range("K54").Value = 9999
$ws.Range("L54").Value = 5000
$ws.Range("M54").Value = -9513
$ws.Range("N54").Value = -5972
$ws.Range("H86").Value = 2494.0715
$ws.Range("I86").Value = 2633.6667
$ws.Range("J86").Value = 2242.8
$ws.Range("K86").Value = 2633.6667
$ws.Range("L86").Value = 2242.8
$ws.Range("M86").Value = -1510.6667
$ws.Range("N86").Value = -4488.8
$ws.Range("H87").Value = 33340.5
$ws.Range("J87").Value = 33340.5
$ws.Range("L87").Value = 33340.5
$ws.Range("N87").Value = -35836.5
$ws.Range("H89").Value = 2494.0715
$ws.Range("I89").Value = 2633.6667
$ws.Range("J89").Value = 2242.8
$ws.Range("K89").Value = 13168.3335
$ws.Range("L89").Value = 11214
$ws.Range("M89").Value = -7552.333500000001
$ws.Range("N89").Value = -22446
$ws.Range("H90").Value = 33340.5
$ws.Range("J90").Value = 33340.5
$ws.Range("L90").Value = 100021.5
$ws.Range("N90").Value = -112501.5
$ws.Range("H98").Value = 1352.1923
$ws.Range("I98").Value = 882.6316
$ws.Range("K98").Value = 882.6316
$ws.Range("M98").Value = 615.3684
$ws.Range("H112").Value = 1493.4375
$ws.Range("J112").Value = 1548.5714
$ws.Range("L112").Value = 4645.7142
$ws.Range("N112").Value = -6861.7142
$ws.Range("H122").Value = 1352.1923
$ws.Range("I122").Value = 882.6316
$ws.Range("K122").Value = 2647.8948
$ws.Range("M122").Value = -197.8948

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1301.0454
$ws.Range("I2").Value = 1077.5238
$ws.Range("K2").Value = 1077.5238
$ws.Range("M2").Value = -964.5237999999999
$ws.Range("H21").Value = 3724.75
$ws.Range("I21").Value = 3238.2
$ws.Range("J21").Value = 4072.2856
$ws.Range("K21").Value = 3238.2
$ws.Range("L21").Value = 4072.2856
$ws.Range("M21").Value = -2864.2
$ws.Range("N21").Value = -4820.2856
$ws.Range("H32").Value = 4288.5225
$ws.Range("I32").Value = 3286.7285
$ws.Range("K32").Value = 3286.7285
$ws.Range("M32").Value = -2999.7285
$ws.Range("H45").Value = 7395.2856
$ws.Range("I45").Value = 13126.667
$ws.Range("K45").Value = 13126.667
$ws.Range("M45").Value = -12749.667
$ws.Range("H61").Value = 1795.6471
$ws.Range("I61").Value = 1737.9667
$ws.Range("J61").Value = 2228.25
$ws.Range("K61").Value = 1737.9667
$ws.Range("L61").Value = 2228.25
$ws.Range("M61").Value = -1525.9667
$ws.Range("N61").Value = -2652.25
$ws.Range("H110").Value = 2257.8484
$ws.Range("I110").Value = 1574.4814
$ws.Range("K110").Value = 1574.4814
$ws.Range("M110").Value = 470.5186000000001
$ws.Range("H116").Value = 1301.0454
$ws.Range("I116").Value = 1077.5238
$ws.Range("K116").Value = 1077.5238
$ws.Range("M116").Value = 1216.4762
$ws.Range("H136").Value = 1795.6471
$ws.Range("I136").Value = 1737.9667
$ws.Range("J136").Value = 2228.25
$ws.Range("K136").Value = 5213.9001
$ws.Range("L136").Value = 6684.75
$ws.Range("M136").Value = -2663.9001
$ws.Range("N136").Value = -11784.75

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1301.0454
$ws.Range("I3").Value = 1077.5238
$ws.Range("K3").Value = 1077.5238
$ws.Range("M3").Value = -963.5237999999999
$ws.Range("H94").Value = 76925816
$ws.Range("I94").Value = 142859360
$ws.Range("J94").Value = 3355
$ws.Range("K94").Value = 142859360
$ws.Range("L94").Value = 3355
$ws.Range("M94").Value = -142858909
$ws.Range("N94").Value = -4257

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 46666830
$ws.Range("I3").Value = 46666830
$ws.Range("K3").Value = 46666830
$ws.Range("M3").Value = -46666717
$ws.Range("H11").Value = 1020000
$ws.Range("J11").Value = 1020000
$ws.Range("L11").Value = 1020000
$ws.Range("N11").Value = -1020280
$ws.Range("H88").Value = 25112.666
$ws.Range("J88").Value = 25112.666
$ws.Range("L88").Value = 25112.666
$ws.Range("N88").Value = -25924.666
$ws.Range("H91").Value = 25112.666
$ws.Range("J91").Value = 25112.666
$ws.Range("L91").Value = 25112.666
$ws.Range("N91").Value = -27920.666
$ws.Range("H132").Value = 2623.5833
$ws.Range("I132").Value = 2499.9412
$ws.Range("K132").Value = 7499.823600000001
$ws.Range("M132").Value = -4969.823600000001
$ws.Range("H141").Value = 234760.8
$ws.Range("J141").Value = 234760.8
$ws.Range("L141").Value = 234760.8
$ws.Range("N141").Value = -245120.8

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1792.6666
$ws.Range("I131").Value = 1403.4
$ws.Range("K131").Value = 4210.200000000001
$ws.Range("M131").Value = 829.7999999999993
$ws.Range("H137").Value = 4169384.2
$ws.Range("J137").Value = 3643.8572
$ws.Range("L137").Value = 10931.5716
$ws.Range("N137").Value = -21131.5716

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H22").Value = 6030
$ws.Range("I22").Value = 10310
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 10310
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -9781
$ws.Range("N22").Value = -2808
$ws.Range("H102").Value = 1731.4
$ws.Range("I102").Value = 1756.6207
$ws.Range("K102").Value = 1756.6207
$ws.Range("M102").Value = -134.6206999999999

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1156.3077
$ws.Range("I55").Value = 1339.7
$ws.Range("J55").Value = 545
$ws.Range("K55").Value = 1339.7
$ws.Range("L55").Value = 545
$ws.Range("M55").Value = -1166.7
$ws.Range("N55").Value = -891
$ws.Range("H132").Value = 4604.8125
$ws.Range("I132").Value = 3675.24
$ws.Range("K132").Value = 11025.72
$ws.Range("M132").Value = -8495.719999999999

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 6500
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H96").Value = 3028
$ws.Range("I96").Value = 1597.5555
$ws.Range("K96").Value = 1597.5555
$ws.Range("M96").Value = -224.5554999999999
$ws.Range("H132").Value = 2952.7114
$ws.Range("I132").Value = 2180.125
$ws.Range("J132").Value = 5528
$ws.Range("K132").Value = 6540.375
$ws.Range("L132").Value = 16584
$ws.Range("M132").Value = -4010.375
$ws.Range("N132").Value = -21644

Write-Host "Applied all updates"